$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "95.884.35"
$ws.Range("E2").Value2 = "  +4.19%  "

Set-TextCell $ws.Range("D3") "3.664.42"
$ws.Range("E3").Value2 = "  +10.15%  "

$ws.Range("E4").Value2 = "  +0.12%  "

Set-TextCell $ws.Range("D5") "242.24"
$ws.Range("E5").Value2 = "  +4.86%  "

Set-TextCell $ws.Range("D6") "645.43"
$ws.Range("E6").Value2 = "  +5.11%  "

$ws.Range("E7").Value2 = "  +4.93%  "

Set-TextCell $ws.Range("D8") "0.401"
$ws.Range("E8").Value2 = "  +3.73%  "

Set-TextCell $ws.Range("D9") "1.00"
$ws.Range("E9").Value2 = "  -0.09%  "

Set-TextCell $ws.Range("D10") "1.01"
$ws.Range("E10").Value2 = "  +5.34%  "

Set-TextCell $ws.Range("D11") "3.663.70"
$ws.Range("E11").Value2 = "  +10.16%  "

Set-TextCell $ws.Range("D12") "43.88"
$ws.Range("E12").Value2 = "  +2.70%  "

$ws.Range("E13").Value2 = "  +3.50%  "

Set-TextCell $ws.Range("D14") "6.38"
$ws.Range("E14").Value2 = "  +3.90%  "

Set-TextCell $ws.Range("D15") "4.365.09"
$ws.Range("E15").Value2 = "  +10.53%  "

Set-TextCell $ws.Range("D16") "95.791.34"
$ws.Range("E16").Value2 = "  +4.27%  "

Set-TextCell $ws.Range("D17") "0.0000257"
$ws.Range("E17").Value2 = "  +5.51%  "

Set-TextCell $ws.Range("D18") "3.667.65"
$ws.Range("E18").Value2 = "  +10.29%  "

Set-TextCell $ws.Range("D19") "13.42"
$ws.Range("E19").Value2 = "  +23.55%  "

Set-TextCell $ws.Range("D20") "8.05"
$ws.Range("E20").Value2 = "  -0.12%  "

Set-TextCell $ws.Range("D21") "18.61"
$ws.Range("E21").Value2 = "  +7.46%  "

Set-TextCell $ws.Range("D22") "519.71"
$ws.Range("E22").Value2 = "  +5.57%  "

Set-TextCell $ws.Range("D23") "3.43"
$ws.Range("E23").Value2 = "  -0.38%  "

Set-TextCell $ws.Range("D24") "0.480"
$ws.Range("E24").Value2 = "  +9.84%  "

$ws.Range("B25").Value2 = "NEARProtocol"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D25") "6.84"
$ws.Range("E25").Value2 = "  +4.54%  "

$ws.Range("B26").Value2 = "PEPE"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D26") "0.0000198"
$ws.Range("E26").Value2 = "  +8.37%  "

Set-TextCell $ws.Range("D27") "93.71"
$ws.Range("E27").Value2 = "  +1.08%  "

Set-TextCell $ws.Range("D28") "12.67"
$ws.Range("E28").Value2 = "  +6.13%  "

Set-TextCell $ws.Range("D29") "3.15"
$ws.Range("E29").Value2 = "  +20.68%  "

Set-TextCell $ws.Range("D30") "11.69"
$ws.Range("E30").Value2 = "  +5.00%  "

$ws.Range("E31").Value2 = "  +2.03%  "

$ws.Range("E32").Value2 = "  -0.09%  "

Set-TextCell $ws.Range("D33") "33.01"
$ws.Range("E33").Value2 = "  +16.53%  "

$ws.Range("E34").Value2 = "  +4.18%  "

Set-TextCell $ws.Range("D35") "1.00"
$ws.Range("E35").Value2 = "  +0.44%  "

Set-TextCell $ws.Range("D36") "0.579"
$ws.Range("E36").Value2 = "  +9.69%  "

Set-TextCell $ws.Range("D37") "561.62"
$ws.Range("E37").Value2 = "  +0.26%  "

Set-TextCell $ws.Range("D38") "7.93"
$ws.Range("E38").Value2 = "  +6.76%  "

$ws.Range("E39").Value2 = "  +9.62%  "

Set-TextCell $ws.Range("D40") "0.967"
$ws.Range("E40").Value2 = "  +11.13%  "

$ws.Range("E41").Value2 = "  +2.86%  "

$ws.Range("E42").Value2 = "  -0.07%  "

$ws.Range("E43").Value2 = "  +4.72%  "

Set-TextCell $ws.Range("D44") "5.81"
$ws.Range("E44").Value2 = "  +7.15%  "

Set-TextCell $ws.Range("D45") "0.0430"
$ws.Range("E45").Value2 = "  +3.65%  "

Set-TextCell $ws.Range("D47") "33.80"
$ws.Range("E47").Value2 = "  +50.37%  "

Set-TextCell $ws.Range("D48") "2.23"
$ws.Range("E48").Value2 = "  +5.18%  "

Set-TextCell $ws.Range("D49") "54.80"
$ws.Range("E49").Value2 = "  +5.29%  "

Set-TextCell $ws.Range("D50") "8.31"
$ws.Range("E50").Value2 = "  +3.74%  "

Set-TextCell $ws.Range("D51") "3.49"
$ws.Range("E51").Value2 = "  -3.59%  "
